$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.394.61'
$ws.Range("E2").Value = '  +0.45%  '
$ws.Range("D3").Value = '3.443.05'
$ws.Range("E3").Value = '  +0.02%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = "'584.22"
$ws.Range("E5").Value = '  +0.76%  '
$ws.Range("D6").Value = "'177.42"
$ws.Range("E6").Value = '  +2.07%  '
$ws.Range("D7").Value = "'0.628"
$ws.Range("E7").Value = '  +6.26%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").Value = '3.442.55'
$ws.Range("E9").Value = '  +0.08%  '
$ws.Range("D10").Value = "'0.133"
$ws.Range("E10").Value = '  +1.09%  '
$ws.Range("E11").Value = '  +1.38%  '
$ws.Range("D12").Value = "'0.416"
$ws.Range("E12").Value = '  -0.34%  '
$ws.Range("D13").Value = '4.049.96'
$ws.Range("E13").Value = '  +0.35%  '
$ws.Range("E14").Value = '  +1.78%  '
$ws.Range("D15").Value = "'29.89"
$ws.Range("E15").Value = '  -0.51%  '
$ws.Range("D16").Value = '66.343.31'
$ws.Range("E16").Value = '  +0.39%  '
$ws.Range("D17").Value = "'0.0000173"
$ws.Range("E17").Value = '  +1.00%  '
$ws.Range("D18").Value = '3.436.52'
$ws.Range("E18").Value = '  -0.23%  '
$ws.Range("D19").Value = "'5.95"
$ws.Range("E19").Value = '  +0.03%  '
$ws.Range("D20").Value = "'13.84"
$ws.Range("E20").Value = '  +0.86%  '
$ws.Range("D21").Value = "'371.35"
$ws.Range("E21").Value = '  -1.39%  '
$ws.Range("D22").Value = "'7.63"
$ws.Range("E22").Value = '  -1.59%  '
$ws.Range("D23").Value = "'73.15"
$ws.Range("E23").Value = '  +1.83%  '
$ws.Range("D24").Value = "'0.0000128"
$ws.Range("E24").Value = '  +8.38%  '
$ws.Range("D25").Value = "'0.997"
$ws.Range("E25").Value = '  -0.18%  '
$ws.Range("D26").Value = "'0.534"
$ws.Range("E26").Value = '  -1.82%  '
$ws.Range("D27").Value = "'9.90"
$ws.Range("E27").Value = '  +1.34%  '
$ws.Range("E28").Value = '  +2.48%  '
$ws.Range("E29").Value = '  +0.19%  '
$ws.Range("D30").Value = "'5.88"
$ws.Range("E30").Value = '  +1.22%  '
$ws.Range("D31").Value = "'1.99"
$ws.Range("E31").Value = '  +0.47%  '
$ws.Range("D32").Value = "'23.60"
$ws.Range("E32").Value = '  -2.00%  '
$ws.Range("E33").Value = '  -0.07%  '
$ws.Range("D34").Value = "'7.07"
$ws.Range("E34").Value = '  -0.34%  '
$ws.Range("E35").Value = '  -3.10%  '
$ws.Range("E36").Value = '  +0.47%  '
$ws.Range("D37").Value = "'163.25"
$ws.Range("E37").Value = '  +2.37%  '
$ws.Range("E38").Value = '  -0.06%  '
$ws.Range("D39").Value = "'27.87"
$ws.Range("E39").Value = '  -4.09%  '
$ws.Range("D40").Value = "'1.80"
$ws.Range("E40").Value = '  +1.97%  '
$ws.Range("D41").Value = "'2.59"
$ws.Range("E41").Value = '  +2.38%  '
$ws.Range("D42").Value = "'4.48"
$ws.Range("E42").Value = '  +0.37%  '
$ws.Range("D43").Value = '2.751.79'
$ws.Range("E43").Value = '  +3.49%  '
$ws.Range("D44").Value = "'6.46"
$ws.Range("E44").Value = '  +1.45%  '
$ws.Range("D45").Value = "'0.0694"
$ws.Range("E45").Value = '  +0.44%  '
$ws.Range("D46").Value = "'25.38"
$ws.Range("E46").Value = '  +4.44%  '
$ws.Range("D47").Value = "'339.41"
$ws.Range("E47").Value = '  +8.72%  '
$ws.Range("D48").Value = "'39.96"
$ws.Range("E48").Value = '  -0.55%  '
$ws.Range("D49").Value = "'0.0286"
$ws.Range("E49").Value = '  -1.02%  '
$ws.Range("D51").Value = "'31.62"
$ws.Range("E51").Value = '  +3.84%  '
